# Append the new data row (row 39) to the ModCounts sheet:
#   A39 = "2025/12/18" (text, same as the other Date cells)
#   B39 = "逃离鸭科夫"  (text, same Game name as every other row)
#   C39 = 1345          (numeric ModCount)
# then copy the formatting (cell style) from the previous data row (38)
# so the new row matches the rest of the table (centered alignment).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use a leading apostrophe so the date-like string is stored as text
# (matching the inlineStr/text representation used by every other row)
# instead of being auto-converted into a date serial number.
$ws.Range("A39").Value = "'2025/12/18"
$ws.Range("B39").Value = "逃离鸭科夫"
$ws.Range("C39").Value = 1345

# Copy the style/formatting of the row above (row 38) onto the new row
# so the new cells end up with the same formatting (e.g. centered text)
# as the rest of the table.
$ws.Range("A38:C38").Copy()
$ws.Range("A39:C39").PasteSpecial(-4122)
